$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C11").Value = "Los Lagos"
$ws.Range("D11").Value = 44503
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100108
$ws.Range("H11").Value = "Tropicales y subtropicales"
$ws.Range("I11").Value = 100108002
$ws.Range("J11").Value = "Mango"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8500
$ws.Range("P11").Value = 8250
$ws.Range("Q11").Value = "$/bandeja 4 kilos"
$ws.Range("R11").Value = "Perú"
$ws.Range("S11").Value = 2062
$ws.Range("T11").Value = 4
